# This workbook lists fantasy-basketball players with their position(s) and team.
# The edit moves the "Deni Avdija" entry (originally on row 16) up to row 7,
# right after "Dillon Brooks" and before "Paolo Banchero". This shifts the rows
# that used to sit between them (Paolo Banchero, Payton Pritchard, Jaylen Brown)
# down by one row. All other rows are unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 7: Deni Avdija (moved up from the old row 16)
$ws.Range("A7").Value = "Deni Avdija"
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "Portland Trail Blazers"

# New row 14: Paolo Banchero (shifted down from the old row 7)
$ws.Range("A14").Value = "Paolo Banchero"
$ws.Range("B14").Value = "SF,PF"
$ws.Range("C14").Value = "Orlando Magic"

# New row 15: Payton Pritchard (shifted down from the old row 14)
$ws.Range("A15").Value = "Payton Pritchard"
$ws.Range("B15").Value = "PG,SG"
$ws.Range("C15").Value = "Boston Celtics"

# New row 16: Jaylen Brown (shifted down from the old row 15)
$ws.Range("A16").Value = "Jaylen Brown"
$ws.Range("B16").Value = "SG,SF"
$ws.Range("C16").Value = "Boston Celtics"
